$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = " Originea musculus extensor digitorum: А. Epicondylus lateralis humeri"
$ws.Range("B4").Value = "B. Epicondylus medialis humeri"
$ws.Range("C4").Value = "C. Olecranon"
$ws.Range("D4").Value = "D. Tuberositas radii E. Fascia antebrachii"
$ws.Range("I4").Value = "Mușchiul extensor al degetelor își are originea pe epicondilul lateral al humerusului și pe fascia antebrațului. Distal mușchiul formează patru tendoane, care se inseră pe fața dorsală a falangelor într-un mod cu totul deosebit.Enunțuri corecte sunt „A” și „E”."

# Row 5
$ws.Range("A5").Value = " Originea musculus flexor carpi radialis: А. Epicondylus lateralis humeri"
$ws.Range("B5").Value = "B. Epicondylus medialis humeri"
$ws.Range("C5").Value = "C. Olecranon"
$ws.Range("D5").Value = "D. Tuberositas radii "
$ws.Range("E5").Value = "E. Fascia antebraţului"
$ws.Range("I5").Value = "Flexorul radial al carpului își ia originea de pe epicondilul medial al humerusului, fascia antebrațului, și de la septele fibroase, care îl separă de pronatorul rotund și palmarul lung și se inseră pe baza osului metacarpian II trecând în prealabil prin canalul radial al carpului de sub retinaculul flexorilor; tendonul lui delimitează șanțul pulsului.Corect – „B” și „E”."

# Row 6
$ws.Range("A6").Value = " Originea musculus flexor carpi radialis: А. Epicondylus lateralis humeri"
$ws.Range("B6").Value = "B. Epicondylus medialis humeri"
$ws.Range("C6").Value = "C. Olecranon"
$ws.Range("D6").Value = "D. Tuberositas radii "
$ws.Range("E6").Value = "E. Fascia antebraţului"
$ws.Range("I6").Value = "Flexorul radial al carpului își ia originea de pe epicondilul medial al humerusului, fascia antebrațului, și de la septele fibroase, care îl separă de pronatorul rotund și palmarul lung și se inseră pe baza osului metacarpian II trecând în prealabil prin canalul radial al carpului de sub retinaculul flexorilor; tendonul lui delimitează șanțul pulsului.Corect – „B” și „E”."

# Row 7
$ws.Range("A7").Value = " Originea musculus flexor carpi radialis: "
$ws.Range("B7").Value = "B. Epicondylus medialis humeri"
$ws.Range("C7").Value = "C. Olecranon"
$ws.Range("D7").Value = "D. Tuberositas radii "
$ws.Range("E7").Value = "E. Fascia antebraţului"
$ws.Range("I7").Value = "А. Epicondylus lateralis humeriА. Epicondylus lateralis humeriFlexorul radial al carpului își ia originea de pe epicondilul medial al humerusului, fascia antebrațului, și de la septele fibroase, care îl separă de pronatorul rotund și palmarul lung și se inseră pe baza osului metacarpian II trecând în prealabil prin canalul radial al carpului de sub retinaculul flexorilor; tendonul lui delimitează șanțul pulsului.Corect – „B” și „E”."

# Row 8
$ws.Range("A8").Value = " Originea musculus flexor carpi radialis: "
$ws.Range("B8").Value = "B. Epicondylus medialis humeri"
$ws.Range("C8").Value = "C. Olecranon"
$ws.Range("D8").Value = "D. Tuberositas radii "
$ws.Range("E8").Value = "E. Fascia antebraţului"
$ws.Range("I8").Value = "А. Epicondylus lateralis humeriА. Epicondylus lateralis humeriFlexorul radial al carpului își ia originea de pe epicondilul medial al humerusului, fascia antebrațului, și de la septele fibroase, care îl separă de pronatorul rotund și palmarul lung și se inseră pe baza osului metacarpian II trecând în prealabil prin canalul radial al carpului de sub retinaculul flexorilor; tendonul lui delimitează șanțul pulsului.Corect – „B” și „E”."

# Row 9
$ws.Range("A9").Value = " Originea musculus flexor carpi radialis: "
$ws.Range("B9").Value = "B. Epicondylus medialis humeri"
$ws.Range("C9").Value = "C. Olecranon"
$ws.Range("D9").Value = "D. Tuberositas radii "
$ws.Range("E9").Value = "E. Fascia antebraţului"
$ws.Range("I9").Value = "А. Epicondylus lateralis humeriА. Epicondylus lateralis humeriFlexorul radial al carpului își ia originea de pe epicondilul medial al humerusului, fascia antebrațului, și de la septele fibroase, care îl separă de pronatorul rotund și palmarul lung și se inseră pe baza osului metacarpian II trecând în prealabil prin canalul radial al carpului de sub retinaculul flexorilor; tendonul lui delimitează șanțul pulsului.Corect – „B” și „E”."

# Row 10
$ws.Range("A10").Value = " Originea musculus flexor carpi radialis: "
$ws.Range("I10").Value = "А. Epicondylus lateralis humeriB. Epicondylus medialis humeriC. OlecranonD. Tuberositas radii E. Fascia antebraţuluiА. Epicondylus lateralis humeriB. Epicondylus medialis humeriC. OlecranonD. Tuberositas radii E. Antebrachial fasciaA. Epicondylus lateralis humeriB. Epicondylus medialis humeriC. OlecranonD. Tuberositas radiiE. Фасция предплечьяFlexorul radial al carpului își ia originea de pe epicondilul medial al humerusului, fascia antebrațului, și de la septele fibroase, care îl separă de pronatorul rotund și palmarul lung și se inseră pe baza osului metacarpian II trecând în prealabil prin canalul radial al carpului de sub retinaculul flexorilor; tendonul lui delimitează șanțul pulsului.Corect – „B” și „E”."

# Row 11
$ws.Range("A11").Value = " Originea musculus flexor carpi radialis: "
$ws.Range("B11").Value = "А. Epicondylus lateralis humeri"
$ws.Range("C11").Value = "B. Epicondylus medialis humeri"
$ws.Range("D11").Value = "C. Olecranon"
$ws.Range("E11").Value = "D. Tuberositas radii "
$ws.Range("F11").Value = "E. Fascia antebraţului"
$ws.Range("I11").Value = "Flexorul radial al carpului își ia originea de pe epicondilul medial al humerusului, fascia antebrațului, și de la septele fibroase, care îl separă de pronatorul rotund și palmarul lung și se inseră pe baza osului metacarpian II trecând în prealabil prin canalul radial al carpului de sub retinaculul flexorilor; tendonul lui delimitează șanțul pulsului.Corect – „B” și „E”."

# Row 12
$ws.Range("A12").Value = " Originea musculus flexor carpi radialis: "
$ws.Range("B12").Value = "А. Epicondylus lateralis humeri"
$ws.Range("C12").Value = "B. Epicondylus medialis humeri"
$ws.Range("D12").Value = "C. Olecranon"
$ws.Range("E12").Value = "D. Tuberositas radii "
$ws.Range("F12").Value = "E. Fascia antebraţului"
$ws.Range("I12").Value = "Flexorul radial al carpului își ia originea de pe epicondilul medial al humerusului, fascia antebrațului, și de la septele fibroase, care îl separă de pronatorul rotund și palmarul lung și se inseră pe baza osului metacarpian II trecând în prealabil prin canalul radial al carpului de sub retinaculul flexorilor; tendonul lui delimitează șanțul pulsului.Corect – „B” și „E”."

# Row 13
$ws.Range("A13").Value = " Originea musculus flexor carpi radialis: "
$ws.Range("B13").Value = "Epicondylus lateralis humeri"
$ws.Range("C13").Value = "Epicondylus medialis humeri"
$ws.Range("D13").Value = "Olecranon"
$ws.Range("E13").Value = "Tuberositas radii"
$ws.Range("F13").Value = "Fascia antebraţului"
$ws.Range("I13").Value = "Flexorul radial al carpului își ia originea de pe epicondilul medial al humerusului, fascia antebrațului, și de la septele fibroase, care îl separă de pronatorul rotund și palmarul lung și se inseră pe baza osului metacarpian II trecând în prealabil prin canalul radial al carpului de sub retinaculul flexorilor; tendonul lui delimitează șanțul pulsului.Corect – „B” și „E”."
